# Commit: "Fruta / hortaliza, semanal"
# The sheet's weekly data table (rows 254..320, one record per row) gets a
# new record inserted at row 255; every existing row from the old 255
# onward shifts down by one (old 255 -> new 256, ..., old 320 -> new 321),
# and the sheet's used-range grows from A1:R320 to A1:R321.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 255, pushing rows 255..320 down to 256..321.
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new record's data.
$ws.Cells.Item(255, 1).Value  = 10
$ws.Cells.Item(255, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(255, 3).Value  = "La Araucanía"
$ws.Cells.Item(255, 4).Value  = 44736
$ws.Cells.Item(255, 5).Value  = 9
$ws.Cells.Item(255, 6).Value  = 100112044
$ws.Cells.Item(255, 7).Value  = "Perejil"
$ws.Cells.Item(255, 8).Value  = "Sin especificar"
$ws.Cells.Item(255, 9).Value  = "Primera"
$ws.Cells.Item(255, 10).Value = 40
$ws.Cells.Item(255, 11).Value = 4000
$ws.Cells.Item(255, 12).Value = 4000
$ws.Cells.Item(255, 13).Value = 4000
$ws.Cells.Item(255, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(255, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(255, 16).Value = 1333
$ws.Cells.Item(255, 17).Value = 3
$ws.Cells.Item(255, 18).Value = "Hortaliza"
